$d = $word.ActiveDocument

# Commit: "FUNCION DE GENERAR CERTIFICADO"
# The placeholder tags {{con_si}} / {{con_no}} are renamed to
# {{conf_si}} / {{conf_no}} (an "f" is inserted after "con").

$d.Content.Find.Execute("{{con_si}}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{conf_si}}", 2)

$d.Content.Find.Execute("{{con_no}}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{conf_no}}", 2)
